$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$srcFormat = $ws.Range("H1")
$srcFormat.Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @(
    @(4, 5),
    @(2, 4),
    @(9, 10),
    @(8, 8),
    @(8, 8),
    @(4, 4),
    @(8, 9),
    @(2, 4),
    @(1, 6),
    @(1, 5),
    @(3, 7),
    @(1, 5),
    @(1, 3),
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(1, 4),
    @(1, 3),
    @(3, 4),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
